$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.794.27'
$ws.Range('E2').Value = '  -1.38%  '
$ws.Range('D3').Value = '1.549.88'
$ws.Range('E3').Value = '  -1.42%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = "'" + '206.03'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.60%  '
$ws.Range('E6').Value = '  -1.64%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  -0.94%  '
$ws.Range('D9').Value = "'" + '21.42'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.58%  '
$ws.Range('E10').Value = '  -1.22%  '
$ws.Range('D11').Value = "'" + '0.0854'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.56%  '
$ws.Range('D12').Value = '1.770.78'
$ws.Range('E12').Value = '  -1.46%  '
$ws.Range('D13').Value = '1.563.80'
$ws.Range('E13').Value = '  -0.35%  '
$ws.Range('E14').Value = '  -2.52%  '
$ws.Range('E15').Value = '  -1.00%  '
$ws.Range('D16').Value = '26.798.74'
$ws.Range('E16').Value = '  -1.45%  '
$ws.Range('D17').Value = "'" + '61.21'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.58%  '
$ws.Range('D18').Value = "'" + '215.04'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.43%  '
$ws.Range('E19').Value = '  +0.37%  '
$ws.Range('E20').Value = '  -1.29%  '
$ws.Range('E21').Value = '  -0.04%  '
$ws.Range('E22').Value = '  -0.59%  '
$ws.Range('E23').Value = '  -4.11%  '
$ws.Range('E24').Value = '  -0.76%  '
$ws.Range('D25').Value = "'" + '153.43'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.51%  '
$ws.Range('D26').Value = "'" + '6.52'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.29%  '
$ws.Range('E27').Value = '  -0.18%  '
$ws.Range('E28').Value = '  -0.05%  '
$ws.Range('E29').Value = '  -1.37%  '
$ws.Range('D30').Value = "'" + '0.0461'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.44%  '
$ws.Range('E31').Value = '  -1.43%  '
$ws.Range('E32').Value = '  +0.24%  '
$ws.Range('D33').Value = '1.345.82'
$ws.Range('E33').Value = '  -4.37%  '
$ws.Range('E34').Value = '  -0.03%  '
$ws.Range('E35').Value = '  -3.10%  '
$ws.Range('E36').Value = '  -0.36%  '
$ws.Range('D37').Value = "'" + '0.933'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.93%  '
$ws.Range('E39').Value = '  +0.90%  '
$ws.Range('E40').Value = '  -1.67%  '
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('E42').Value = '  +5.25%  '
$ws.Range('E43').Value = '  -0.03%  '
$ws.Range('E44').Value = '  +0.32%  '
$ws.Range('E45').Value = '  -3.43%  '
$ws.Range('D46').Value = "'" + '62.98'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.06%  '
$ws.Range('D47').Value = '1.685.00'
$ws.Range('E47').Value = '  -1.41%  '
$ws.Range('E48').Value = '  -2.87%  '
$ws.Range('D49').Value = "'" + '85.81'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.01%  '
$ws.Range('D50').Value = "'" + '0.0511'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.42%  '
$ws.Range('D51').Value = '0.0₇0972'
$ws.Range('E51').Value = '  -0.39%  '
